$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ewrd"
$ws.Range("B1").Value = "fewfr"
$ws.Range("C1").Value = "fedcf"
$ws.Range("D1").Value = 23
$ws.Range("E1").Value = "d"

$ws.Range("A2").Value = "ewrd"
$ws.Range("B2").Value = "fewfr"
$ws.Range("C2").Value = "fedcf"
$ws.Range("D2").Value = 23
$ws.Range("E2").Value = "d"
